$d = $word.ActiveDocument

$d.Content.Find.Execute("This paragraph centered", $true, $false, $false, $false, $false,
                         $true, 1, $false, "This paragraph centered", 2)
